$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-0.12839389494417475, 0.12833192244811897),
  @(-0.12595429452686169, 0.1257056409456716),
  @(-0.076003000305535195, 0.075812318505009912),
  @(-0.067812318527879611, 0.067359810333744718),
  @(-0.064359810345857582, 0.062822722351545401),
  @(-0.012057899823764373, 0.011908558314157602),
  @(-0.0019085583455806798, 0.0018844350237183782),
  @(0.0081155649447532952, -0.0081374785888814749),
  @(-0.027808142841815098, 0.027588325231679978),
  @(-0.025588325245371024, 0.025573093311754391),
  @(-0.024386016686771761, 0.024359802617818893),
  @(-0.020859802635417424, 0.020668741812290392),
  @(-0.017168741831359746, 0.017081513816760818),
  @(-0.0090815138468984813, 0.0090528277403523205),
  @(-0.0080528277541924709, 0.0080342871564313967),
  @(-0.0060342871728966685, 0.0060036415158770673),
  @(-0.0040036415327353581, 0.0039999999783146833),
  @(0.0036223720642176716, -0.0036590741938660187),
  @(0.0076590741840711907, -0.0079035286184732101),
  @(-0.0080171183275314917, 0.0080057044768864216),
  @(-0.0040057044865395852, 0.0039999999902580186),
  @(-0.045702640652486082, 0.045491924144830165),
  @(-0.040491924160342307, 0.040097689919703861),
  @(-0.020097689972140564, 0.019999999946863412),
  @(-0.066531039155973914, 0.066443966539896593),
  @(-0.063943966553912546, 0.063834721638137637),
  @(-0.061334721652834823, 0.060706469413871478),
  @(-0.058706469430593877, 0.05828944539173353),
  @(-0.051289445422262325, 0.05117748755439866),
  @(-0.021169058663374596, 0.021023824447650608),
  @(-0.014023824480274172, 0.014001446289830355),
  @(-0.004001446329729319, 0.003999999974221069)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $ws.Cells.Item($i + 1, 1).Value = $data[$i][0]
  $ws.Cells.Item($i + 1, 2).Value = $data[$i][1]
}

$ws.Columns.Item(1).ColumnWidth = 14.625
$ws.Columns.Item(2).ColumnWidth = 14.625
